$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.100786805152893
$ws.Range("B1").Value = 2.207854270935059
$ws.Range("C1").Value = 9.994098663330078
$ws.Range("D1").Value = 1.253017544746399
$ws.Range("E1").Value = 1.255845427513123
